# PT's 5000m performance from 1996
# Adds a new data row (row 13) to the "Other" sheet with Pete Thompson's
# 1996 5000m performance, and updates the sheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry: Performance | Date | Name | (blank) | Po10 Event | Fixture | (blank) | Gender | Age Code | Notes
$ws.Range("B13").Value = "15:30"
$ws.Range("C13").Value = "3 Aug 1996"
$ws.Range("D13").Value = "Pete Thompson"
$ws.Range("F13").Value = "5000"
$ws.Range("G13").Value = "Portsmouth"
$ws.Range("I13").Value = "M"
$ws.Range("J13").Value = "V35"
$ws.Range("K13").Value = "From Noel Moss 9Apr2024"

# Move the active selection to reflect where the editor ended up.
$ws.Range("K14").Select()
